$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1929.1111
$ws.Cells.Item(28, 9).Value = 1545.375
$ws.Cells.Item(28, 10).Value = 4999
$ws.Cells.Item(28, 11).Value = 1545.375
$ws.Cells.Item(28, 12).Value = 4999
$ws.Cells.Item(28, 13).Value = -1060.375
$ws.Cells.Item(28, 14).Value = -5969
$ws.Cells.Item(112, 8).Value = 2579.2104
$ws.Cells.Item(112, 9).Value = 3078
$ws.Cells.Item(112, 10).Value = 2401.0715
$ws.Cells.Item(112, 11).Value = 9234
$ws.Cells.Item(112, 12).Value = 7203.2145
$ws.Cells.Item(112, 13).Value = -8126
$ws.Cells.Item(112, 14).Value = -9419.2145
$ws.Cells.Item(136, 8).Value = 67676
$ws.Cells.Item(136, 10).Value = 67676
$ws.Cells.Item(136, 12).Value = 67676
$ws.Cells.Item(136, 14).Value = -77876
$ws.Cells.Item(137, 8).Value = 8485.375
$ws.Cells.Item(137, 10).Value = 2381
$ws.Cells.Item(137, 12).Value = 7143
$ws.Cells.Item(137, 14).Value = -12243
$ws.Cells.Item(138, 8).Value = 367144.8
$ws.Cells.Item(138, 9).Value = 549048.75
$ws.Cells.Item(138, 10).Value = 3336.9167
$ws.Cells.Item(138, 11).Value = 1647146.25
$ws.Cells.Item(138, 12).Value = 10010.7501
$ws.Cells.Item(138, 13).Value = -1642006.25
$ws.Cells.Item(138, 14).Value = -20290.7501
$ws.Cells.Item(141, 8).Value = 6284.5557
$ws.Cells.Item(141, 9).Value = 6070.125
$ws.Cells.Item(141, 10).Value = 8000
$ws.Cells.Item(141, 11).Value = 18210.375
$ws.Cells.Item(141, 12).Value = 24000
$ws.Cells.Item(141, 13).Value = -13030.375
$ws.Cells.Item(141, 14).Value = -34360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 167813
$ws.Cells.Item(2, 10).Value = 335133
$ws.Cells.Item(2, 12).Value = 335133
$ws.Cells.Item(2, 14).Value = -335359
$ws.Cells.Item(32, 8).Value = 7309.3037
$ws.Cells.Item(32, 9).Value = 7286.54
$ws.Cells.Item(32, 11).Value = 7286.54
$ws.Cells.Item(32, 13).Value = -6999.54
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 8193.742
$ws.Cells.Item(61, 9).Value = 9364.772000000001
$ws.Cells.Item(61, 11).Value = 9364.772000000001
$ws.Cells.Item(61, 13).Value = -9152.772000000001
$ws.Cells.Item(74, 8).Value = 4383.5
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 4383.5
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 4383.5
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 14).Value = -6131.5
$ws.Cells.Item(77, 8).Value = 4383.5
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 4383.5
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 21917.5
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(77, 14).Value = -30653.5
$ws.Cells.Item(116, 8).Value = 167813
$ws.Cells.Item(116, 10).Value = 335133
$ws.Cells.Item(116, 12).Value = 335133
$ws.Cells.Item(116, 14).Value = -339721
$ws.Cells.Item(122, 8).Value = 913345.7
$ws.Cells.Item(122, 9).Value = 3685.6086
$ws.Cells.Item(122, 10).Value = 3005564
$ws.Cells.Item(122, 11).Value = 11056.8258
$ws.Cells.Item(122, 12).Value = 9016692
$ws.Cells.Item(122, 13).Value = -8606.825800000001
$ws.Cells.Item(122, 14).Value = -9021592
$ws.Cells.Item(132, 8).Value = 2529.5818
$ws.Cells.Item(132, 9).Value = 2186.1025
$ws.Cells.Item(132, 10).Value = 3366.8125
$ws.Cells.Item(132, 11).Value = 6558.3075
$ws.Cells.Item(132, 12).Value = 10100.4375
$ws.Cells.Item(132, 13).Value = -4028.3075
$ws.Cells.Item(132, 14).Value = -15160.4375
$ws.Cells.Item(136, 8).Value = 8193.742
$ws.Cells.Item(136, 9).Value = 9364.772000000001
$ws.Cells.Item(136, 11).Value = 28094.316
$ws.Cells.Item(136, 13).Value = -25544.316
$ws.Cells.Item(138, 8).Value = 55203
$ws.Cells.Item(138, 10).Value = 55203
$ws.Cells.Item(138, 12).Value = 55203
$ws.Cells.Item(138, 14).Value = -65483

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 167813
$ws.Cells.Item(3, 10).Value = 335133
$ws.Cells.Item(3, 12).Value = 335133
$ws.Cells.Item(3, 14).Value = -335361

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 14702.357
$ws.Cells.Item(7, 9).Value = 20308.3
$ws.Cells.Item(7, 11).Value = 20308.3
$ws.Cells.Item(7, 13).Value = -20195.3
$ws.Cells.Item(31, 8).Value = 6807.353
$ws.Cells.Item(31, 9).Value = 7577.0454
$ws.Cells.Item(31, 10).Value = 5396.25
$ws.Cells.Item(31, 11).Value = 7577.0454
$ws.Cells.Item(31, 12).Value = 5396.25
$ws.Cells.Item(31, 13).Value = -7282.0454
$ws.Cells.Item(31, 14).Value = -5986.25
$ws.Cells.Item(34, 8).Value = 6807.353
$ws.Cells.Item(34, 9).Value = 7577.0454
$ws.Cells.Item(34, 10).Value = 5396.25
$ws.Cells.Item(34, 11).Value = 7577.0454
$ws.Cells.Item(34, 12).Value = 5396.25
$ws.Cells.Item(34, 13).Value = -7375.0454
$ws.Cells.Item(34, 14).Value = -5800.25
$ws.Cells.Item(122, 8).Value = 10057.869
$ws.Cells.Item(122, 9).Value = 9255
$ws.Cells.Item(122, 11).Value = 27765
$ws.Cells.Item(122, 13).Value = -25315

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 59633116
$ws.Cells.Item(4, 9).Value = 56321150
$ws.Cells.Item(4, 10).Value = 70120990
$ws.Cells.Item(4, 11).Value = 168963450
$ws.Cells.Item(4, 12).Value = 210362970
$ws.Cells.Item(4, 13).Value = -168963338
$ws.Cells.Item(4, 14).Value = -210363194
$ws.Cells.Item(12, 8).Value = 99.411766
$ws.Cells.Item(12, 9).Value = 182
$ws.Cells.Item(12, 10).Value = 26
$ws.Cells.Item(12, 11).Value = 546
$ws.Cells.Item(12, 12).Value = 78
$ws.Cells.Item(12, 13).Value = -373
$ws.Cells.Item(12, 14).Value = -424
$ws.Cells.Item(23, 8).Value = 630.1
$ws.Cells.Item(23, 9).Value = 660
$ws.Cells.Item(23, 10).Value = 610.1667
$ws.Cells.Item(23, 11).Value = 1980
$ws.Cells.Item(23, 12).Value = 1830.5001
$ws.Cells.Item(23, 13).Value = -1745
$ws.Cells.Item(23, 14).Value = -2300.5001
$ws.Cells.Item(86, 8).Value = 413.5
$ws.Cells.Item(86, 10).Value = 419.75
$ws.Cells.Item(86, 12).Value = 1259.25
$ws.Cells.Item(86, 14).Value = -3631.25
$ws.Cells.Item(89, 8).Value = 413.5
$ws.Cells.Item(89, 10).Value = 419.75
$ws.Cells.Item(89, 12).Value = 3777.75
$ws.Cells.Item(89, 14).Value = -15633.75
$ws.Cells.Item(129, 8).Value = 2462.647
$ws.Cells.Item(129, 10).Value = 2993.4614
$ws.Cells.Item(129, 12).Value = 8980.3842
$ws.Cells.Item(129, 14).Value = -18980.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 12072.182
$ws.Cells.Item(113, 9).Value = 15235.125
$ws.Cells.Item(113, 10).Value = 3637.6667
$ws.Cells.Item(113, 11).Value = 15235.125
$ws.Cells.Item(113, 12).Value = 3637.6667
$ws.Cells.Item(113, 13).Value = -13065.125
$ws.Cells.Item(113, 14).Value = -7977.6667
$ws.Cells.Item(134, 8).Value = 91859.60000000001
$ws.Cells.Item(134, 10).Value = 91859.60000000001
$ws.Cells.Item(134, 12).Value = 275578.8
$ws.Cells.Item(134, 14).Value = -280648.8
$ws.Cells.Item(141, 8).Value = 78997.17999999999
$ws.Cells.Item(141, 10).Value = 78997.17999999999
$ws.Cells.Item(141, 12).Value = 78997.17999999999
$ws.Cells.Item(141, 14).Value = -89357.17999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1664
$ws.Cells.Item(113, 9).Value = 760.7692
$ws.Cells.Item(113, 10).Value = 2731.4546
$ws.Cells.Item(113, 11).Value = 2282.3076
$ws.Cells.Item(113, 12).Value = 8194.363799999999
$ws.Cells.Item(113, 13).Value = -112.3076000000001
$ws.Cells.Item(113, 14).Value = -12534.3638
$ws.Cells.Item(141, 8).Value = 104102.29
$ws.Cells.Item(141, 10).Value = 96619.336
$ws.Cells.Item(141, 12).Value = 96619.336
$ws.Cells.Item(141, 14).Value = -106979.336
